$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary: "affecting all Black and Asian-American voters" ->
#    "affecting 50M voters"  (stays a single run, no formatting change)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: split the run so "50M" becomes its own bold,
#    colored run - "... affecting all Black and Asian-American voters,
#    developed ..." -> "... affecting " + bold("50M") + " voters, developed ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.Find.Execute("race coding errors affecting ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.Find.Execute("all Black and Asian-American", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "50M"
$r.Font.Bold = $true
$r.Font.Color = 5258796   # 0x2C3E50 (BGR-packed wdColor for RGB 2C3E50)

# ---------------------------------------------------------------------------
# 3) Reorder work-experience blocks: move the "Analytics Supervisor - GSD&M"
#    block (heading + 4 paragraphs) so it follows "Data Products Manager -
#    Helm/Murmuration" instead of preceding it.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$blockStartPara = $null
$blockEndPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "Analytics Supervisor - GSD*") {
        $blockStartPara = $i
    }
    if ($blockStartPara -ne $null -and $i -ge $blockStartPara -and $t -like "*Advanced Statistical and ML techniques*") {
        $blockEndPara = $i
        break
    }
}

$first = $paras.Item($blockStartPara)
$last = $paras.Item($blockEndPara)
$moveRange = $d.Range($first.Range.Start, $last.Range.End)
$moveRange.Cut() | Out-Null

# Find the destination: right before the "Senior Analyst - Myers Research" heading
$paras = $d.Paragraphs
$destIndex = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "Senior Analyst - Myers Research*") {
        $destIndex = $i
        break
    }
}
$dest = $paras.Item($destIndex)
$insertionRange = $d.Range($dest.Range.Start, $dest.Range.Start)
$insertionRange.Paste() | Out-Null

# Pasting can drop the paragraph-mark style (Heading3) on the re-inserted
# heading paragraph; restore it explicitly.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Analytics Supervisor - GSD*") {
        $p.Style = "Heading 3"
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Key Projects "Impact" line: "affecting all Black and Asian-American
#    voters, improved" -> "affecting 50M voters nationwide, improved"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2) | Out-Null
